$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.858.91'
$ws.Range("E2").Value = '  +0.37%  '

$ws.Range("D3").Value = '1.632.63'
$ws.Range("E3").Value = '  +0.27%  '

$ws.Range("E4").Value = '  +0.54%  '

$ws.Range("D5").Value = '214.44'
$ws.Range("E5").Value = '  +0.33%  '

$ws.Range("E6").Value = '  +1.06%  '

$ws.Range("E7").Value = '  +0.50%  '

$ws.Range("E8").Value = '  -0.49%  '

$ws.Range("D9").Value = '0.0633'
$ws.Range("E9").Value = '  +0.58%  '

$ws.Range("D10").Value = '19.55'
$ws.Range("E10").Value = '  -0.12%  '

$ws.Range("D11").Value = '0.0792'
$ws.Range("E11").Value = '  +0.53%  '

$ws.Range("B12").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C12").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D12").Value = '1.857.32'
$ws.Range("E12").Value = '  +0.19%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '4.24'
$ws.Range("E13").Value = '  -0.15%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.640.11'
$ws.Range("E14").Value = '  +2.43%  '

$ws.Range("D15").Value = '0.544'
$ws.Range("E15").Value = '  -1.02%  '

$ws.Range("D16").Value = '0.0₃0756'
$ws.Range("E16").Value = '  +0.03%  '

$ws.Range("D17").Value = '62.54'
$ws.Range("E17").Value = '  -0.16%  '

$ws.Range("D18").Value = '25.859.45'
$ws.Range("E18").Value = '  +0.29%  '

$ws.Range("E19").Value = '  +0.51%  '

$ws.Range("D20").Value = '193.60'
$ws.Range("E20").Value = '  +1.56%  '

$ws.Range("E21").Value = '  -0.63%  '

$ws.Range("D22").Value = '9.92'
$ws.Range("E22").Value = '  +0.00%  '

$ws.Range("D23").Value = '6.27'
$ws.Range("E23").Value = '  -0.02%  '

$ws.Range("E24").Value = '  +1.34%  '

$ws.Range("D25").Value = '143.43'
$ws.Range("E25").Value = '  +0.93%  '

$ws.Range("E26").Value = '  +0.56%  '

$ws.Range("E27").Value = '  +2.66%  '

$ws.Range("D28").Value = '6.84'
$ws.Range("E28").Value = '  +0.27%  '

$ws.Range("D29").Value = '15.42'
$ws.Range("E29").Value = '  -0.62%  '

$ws.Range("E30").Value = '  +0.34%  '

$ws.Range("E31").Value = '  +1.16%  '

$ws.Range("E32").Value = '  -0.36%  '

$ws.Range("E33").Value = '  +0.02%  '

$ws.Range("E34").Value = '  -1.95%  '

$ws.Range("E35").Value = '  +1.60%  '

$ws.Range("D36").Value = '0.901'
$ws.Range("E36").Value = '  -0.06%  '

$ws.Range("D37").Value = '1.138.13'
$ws.Range("E37").Value = '  -0.18%  '

$ws.Range("E38").Value = '  +0.38%  '

$ws.Range("E39").Value = '  -0.78%  '

$ws.Range("E40").Value = '  +0.49%  '

$ws.Range("E41").Value = '  +0.61%  '

$ws.Range("D42").Value = '99.24'
$ws.Range("E42").Value = '  -1.09%  '

$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").Value = '0.798'
$ws.Range("E43").Value = '  -0.07%  '

$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = '5.42'
$ws.Range("E44").Value = '  -3.21%  '

$ws.Range("D45").Value = '1.766.52'
$ws.Range("E45").Value = '  +0.11%  '

$ws.Range("E46").Value = '  -0.04%  '

$ws.Range("D47").Value = '56.23'
$ws.Range("E47").Value = '  +1.69%  '

$ws.Range("E48").Value = '  +3.32%  '

$ws.Range("E49").Value = '  -0.77%  '

$ws.Range("E50").Value = '  -0.11%  '

$ws.Range("D51").Value = '7.63'
$ws.Range("E51").Value = '  +1.38%  '
